# Generate Report for Archive
#
# The localization run moved the "Overview"/"zh-cn"/"de-de" sheets from
# "Ready for handoff" to "In Translation", and the Status-ish columns that
# used to hold that (now shorter) text were re-sized to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the columns that held the status text to fit the shorter string.
$wsOverview.Range("E1").ColumnWidth = 13.4101845877511
$wsOverview.Range("F1").ColumnWidth = 13.4101845877511
$wsZhCn.Range("C1").ColumnWidth = 13.4101845877511
$wsDeDe.Range("C1").ColumnWidth = 13.4101845877511
